$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F6").Value = 76
$ws.Range("G6").Value = 2270.88
$ws.Range("B10").Value = 28285.77
$ws.Range("F68").Value = 47
$ws.Range("G68").Value = 5410.64
$ws.Range("F77").Value = 257
$ws.Range("G77").Value = 12012.18
$ws.Range("F78").Value = 40
$ws.Range("G78").Value = 2276
$ws.Range("F81").Value = 11
$ws.Range("G81").Value = 336.38
$ws.Range("F84").Value = 33
$ws.Range("G84").Value = 3381.18
$ws.Range("F86").Value = 64
$ws.Range("G86").Value = 8030.08
$ws.Range("B90").Value = 181004.59
$ws.Range("F115").Value = 203
$ws.Range("G115").Value = 19652.43
$ws.Range("B117").Value = 13510.31
$ws.Range("F144").Value = 1050
$ws.Range("G144").Value = 8872.5
$ws.Range("F145").Value = 472
$ws.Range("G145").Value = 3771.28
$ws.Range("F146").Value = 23
$ws.Range("G146").Value = 1936.37
$ws.Range("B147").Value = 14580.15
$ws.Range("F149").Value = 230
$ws.Range("G149").Value = 14904
$ws.Range("B156").Value = 31630.59
$ws.Range("B192").Value = 48706
$ws.Range("E192").Value = 39.8
$ws.Range("F192").Value = -144
$ws.Range("G192").Value = -4795.2
$ws.Range("B193").Value = 64973
$ws.Range("E193").Value = 35.4
$ws.Range("F193").Value = 2
$ws.Range("G193").Value = 66.59999999999999
$ws.Range("F197").Value = 20
$ws.Range("G197").Value = 1241.2
$ws.Range("B216").Value = 40282.97
$ws.Range("B243").Value = 60325
$ws.Range("E243").Value = 151.57
$ws.Range("F243").Value = -102
$ws.Range("G243").Value = -12939.72
$ws.Range("B244").Value = 63560
$ws.Range("E244").Value = 134.87
$ws.Range("F244").Value = 1
$ws.Range("G244").Value = 126.86
$ws.Range("F247").Value = 142
$ws.Range("G247").Value = 14755.22
$ws.Range("F249").Value = 139
$ws.Range("G249").Value = 19156.98
$ws.Range("F252").Value = 2
$ws.Range("G252").Value = 42.06
$ws.Range("F255").Value = 563
$ws.Range("G255").Value = 96458.78999999999
$ws.Range("B260").Value = 190409.13
$ws.Range("F293").Value = 38
$ws.Range("G293").Value = 2672.16
$ws.Range("F300").Value = 164
$ws.Range("G300").Value = 20485.24
$ws.Range("F302").Value = 45
$ws.Range("G302").Value = 9490.049999999999
$ws.Range("F303").Value = 30
$ws.Range("G303").Value = 6326.7
$ws.Range("B304").Value = 175071.7
$ws.Range("B322").Value = 58047
$ws.Range("D322").Value = 105.54
$ws.Range("E322").Value = 126.1
$ws.Range("F322").Value = 39
$ws.Range("G322").Value = 4116.06
$ws.Range("B323").Value = 47097
$ws.Range("D323").Value = 112.28
$ws.Range("E323").Value = 134.16
$ws.Range("F323").Value = 15
$ws.Range("G323").Value = 1684.2
$ws.Range("F334").Value = 193
$ws.Range("G334").Value = 10001.26
$ws.Range("F342").Value = 140
$ws.Range("G342").Value = 4433.8
$ws.Range("F345").Value = 55
$ws.Range("G345").Value = 3377.55
$ws.Range("B346").Value = 26163.7
$ws.Range("F350").Value = 63
$ws.Range("G350").Value = 4833.99
$ws.Range("B358").Value = 35555.12
$ws.Range("B364").Value = 53602
$ws.Range("E364").Value = 15.69
$ws.Range("F364").Value = -231
$ws.Range("G364").Value = -3037.65
$ws.Range("B365").Value = 65068
$ws.Range("E365").Value = 13.97
$ws.Range("F365").Value = 63
$ws.Range("G365").Value = 828.45
$ws.Range("B366").Value = 65066
$ws.Range("E366").Value = 13.61
$ws.Range("F366").Value = 90
$ws.Range("G366").Value = 1152.9
$ws.Range("B367").Value = 53263
$ws.Range("E367").Value = 15.29
$ws.Range("F367").Value = -309
$ws.Range("G367").Value = -3958.29
$ws.Range("B372").Value = 45706
$ws.Range("E372").Value = 23.58
$ws.Range("F372").Value = -202
$ws.Range("G372").Value = -3985.46
$ws.Range("B373").Value = 64922
$ws.Range("E373").Value = 20.98
$ws.Range("F373").Value = 67
$ws.Range("G373").Value = 1321.91
$ws.Range("B375").Value = 45718
$ws.Range("E375").Value = 19.38
$ws.Range("F375").Value = -294
$ws.Range("G375").Value = -4768.68
$ws.Range("B376").Value = 64927
$ws.Range("E376").Value = 17.26
$ws.Range("F376").Value = 106
$ws.Range("G376").Value = 1719.32
$ws.Range("B380").Value = 64925
$ws.Range("E380").Value = 13.97
$ws.Range("F380").Value = 111
$ws.Range("G380").Value = 1459.65
$ws.Range("B381").Value = 45709
$ws.Range("E381").Value = 15.69
$ws.Range("F381").Value = -300
$ws.Range("G381").Value = -3945
$ws.Range("B382").Value = 45702
$ws.Range("E382").Value = 31.43
$ws.Range("F382").Value = -215
$ws.Range("G382").Value = -5654.5
$ws.Range("B383").Value = 64919
$ws.Range("E383").Value = 27.97
$ws.Range("F383").Value = 61
$ws.Range("G383").Value = 1604.3
$ws.Range("B385").Value = 65067
$ws.Range("E385").Value = 15.65
$ws.Range("F385").Value = 126
$ws.Range("G385").Value = 1855.98
$ws.Range("B386").Value = 53595
$ws.Range("E386").Value = 17.61
$ws.Range("F386").Value = -335
$ws.Range("G386").Value = -4934.55
$ws.Range("F400").Value = 0
$ws.Range("G400").Value = 0
$ws.Range("B411").Value = 7628.12
$ws.Range("F434").Value = 9
$ws.Range("G434").Value = 293.76
$ws.Range("B435").Value = 439.24
$ws.Range("B442").Value = 64810
$ws.Range("E442").Value = 291.22
$ws.Range("F442").Value = 4
$ws.Range("G442").Value = 1095.68
$ws.Range("B443").Value = 53319
$ws.Range("E443").Value = 310.64
$ws.Range("F443").Value = -6
$ws.Range("G443").Value = -1643.52
$ws.Range("F471").Value = 29
$ws.Range("G471").Value = 2437.74
$ws.Range("B473").Value = 64830
$ws.Range("E473").Value = 34.9
$ws.Range("F473").Value = 107
$ws.Range("G473").Value = 3512.81
$ws.Range("B474").Value = 60022
$ws.Range("E474").Value = 37.22
$ws.Range("F474").Value = -113
$ws.Range("G474").Value = -3709.79
$ws.Range("B475").Value = 45353.88
$ws.Range("F480").Value = 120
$ws.Range("G480").Value = 14229.6
$ws.Range("B488").Value = 30626.94
$ws.Range("F491").Value = 22
$ws.Range("G491").Value = 3916.44
$ws.Range("B493").Value = 12141.09
$ws.Range("F508").Value = 58
$ws.Range("G508").Value = 6028.52
$ws.Range("B510").Value = 23631.74
$ws.Range("B572").Value = 65079
$ws.Range("F572").Value = 6
$ws.Range("G572").Value = 245.22
$ws.Range("B573").Value = 65362
$ws.Range("F573").Value = 20
$ws.Range("G573").Value = 817.4
$ws.Range("F577").Value = 62
$ws.Range("G577").Value = 2665.38
$ws.Range("F580").Value = 57
$ws.Range("G580").Value = 3248.43
$ws.Range("B583").Value = 16689.95
$ws.Range("F599").Value = 1639
$ws.Range("G599").Value = 267337.29
$ws.Range("F601").Value = 410
$ws.Range("G601").Value = 115976.7
$ws.Range("F602").Value = 332
$ws.Range("G602").Value = 48023.8
$ws.Range("B606").Value = 432185.84
$ws.Range("B619").Value = 1751231.2
$ws.Range("B620").Value = 1751231.2
